$wb = $excel.ActiveWorkbook

# Duplicate an existing parameter sheet so the new sheet inherits the same
# sheet-level formatting (sheetPr/outlinePr, sheetFormatPr, per-cell style
# "1", and the worksheet -> drawing relationship) used throughout this
# workbook, then strip its data and rebuild it with the screening-strategy
# parameters.
$template = $wb.Worksheets.Item("strategy_params")
$template.Copy($null, $template)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "screening_strategies"

# Drop the copied values but keep the per-cell formatting (style index)
# intact so the new sheet's cells still carry s="1" like the rest of the
# workbook.
$ws.Cells.ClearContents()

# Column widths: col A and col C are custom-sized on this sheet; column B
# keeps the default width (unlike the template sheet), so reset the
# inherited column widths first.
$ws.Columns.Item(1).ColumnWidth = 15.29
$ws.Columns.Item(2).ColumnWidth = $template.Columns.Item(2).ColumnWidth
$ws.Columns.Item(3).ColumnWidth = 30.17
$ws.Range("B:B").EntireColumn.AutoFit() | Out-Null

# Header row
$ws.Range("A1").Value = "ParameterName"
$ws.Range("B1").Value = "BaseValue"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Source"

# Data rows
$data = @(
    @("screening_3_cost", 24.0, "Cost of 3 Pap tests ($8 each)"),
    @("screening_10_cost", 80.0, "Cost of 10 Pap tests ($8 each)"),
    @("screening_3_utility", 0.99, "Utility during Pap test year (3 tests)"),
    @("screening_10_utility", 0.99, "Utility during Pap test year (10 tests)"),
    @("screening_age_start", 18.0, "Age when screening begins"),
    @("screening_3_freq", 3.0, "Number of screenings in 3-test strategy"),
    @("screening_10_freq", 10.0, "Number of screenings in 10-test strategy")
)

$row = 2
foreach ($item in $data) {
    # Bring over the header row's formatting (style "1") one row at a time
    # so every new row of cells matches the template's look, including for
    # row 8 which falls outside the range the template sheet originally had
    # formatted.
    $ws.Range("A1:D1").Copy() | Out-Null
    $ws.Range("A$row`:D$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$row").Value = $item[0]
    $ws.Range("B$row").Value = $item[1]
    $ws.Range("C$row").Value = $item[2]
    $row++
}

$excel.CutCopyMode = $false
